$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026535994458571
$ws.Range("D2").Value = 1.030845768367896
$ws.Range("E2").Value = 1.030165566476787
$ws.Range("F2").Value = 1.037166872292359
$ws.Range("I2").Value = 1.034037997230285
$ws.Range("J2").Value = 1.031698890435828
$ws.Range("K2").Value = 1.033655599648712
$ws.Range("L2").Value = 1.032977367706252
$ws.Range("M2").Value = 1.039958531108289
$ws.Range("N2").Value = 1.014589358603693
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027349334693291
$ws.Range("D3").Value = 1.031436204382928
$ws.Range("E3").Value = 1.030926108454984
$ws.Range("F3").Value = 1.038242059231663
$ws.Range("I3").Value = 1.03421903866973
$ws.Range("J3").Value = 1.032152733101873
$ws.Range("K3").Value = 1.034055061067428
$ws.Range("L3").Value = 1.033546335848341
$ws.Range("M3").Value = 1.040842766802175
$ws.Range("N3").Value = 1.014741129812803
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.02787617975084
$ws.Range("D4").Value = 1.031818706837025
$ws.Range("E4").Value = 1.031419132226238
$ws.Range("F4").Value = 1.03893871727807
$ws.Range("I4").Value = 1.034335279573637
$ws.Range("J4").Value = 1.032446300610625
$ws.Range("K4").Value = 1.034313291789549
$ws.Range("L4").Value = 1.033914736935971
$ws.Range("M4").Value = 1.041415272907436
$ws.Range("N4").Value = 1.014839262753253
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028097797874449
$ws.Range("D5").Value = 1.031979617203184
$ws.Range("E5").Value = 1.031626613713016
$ws.Range("F5").Value = 1.03923181645785
$ws.Range("I5").Value = 1.034383929981603
$ws.Range("J5").Value = 1.032569691596269
$ws.Range("K5").Value = 1.034421791650694
$ws.Range("L5").Value = 1.034069669034797
$ws.Range("M5").Value = 1.041656036485401
$ws.Range("N5").Value = 1.014880499913449
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028135016256814
$ws.Range("D6").Value = 1.032006640950493
$ws.Range("E6").Value = 1.031661463284191
$ws.Range("F6").Value = 1.039281042215909
$ws.Range("I6").Value = 1.034392085840891
$ws.Range("J6").Value = 1.032590408000214
$ws.Range("K6").Value = 1.034440005685481
$ws.Range("L6").Value = 1.034095686077595
$ws.Range("M6").Value = 1.041696466531605
$ws.Range("N6").Value = 1.014887422750915
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027879140504438
$ws.Range("D7").Value = 1.031820856513413
$ws.Range("E7").Value = 1.031421903764139
$ws.Range("F7").Value = 1.038942632804126
$ws.Range("I7").Value = 1.034335930497238
$ws.Range("J7").Value = 1.032447949465801
$ws.Range("K7").Value = 1.034314741808428
$ws.Range("L7").Value = 1.033916806926782
$ws.Range("M7").Value = 1.041418489681066
$ws.Range("N7").Value = 1.014839813837163
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026810749843939
$ws.Range("D8").Value = 1.031045214663455
$ws.Range("E8").Value = 1.03042240735422
$ws.Range("F8").Value = 1.037530041505022
$ws.Range("I8").Value = 1.034099368056436
$ws.Range("J8").Value = 1.031852288545628
$ws.Range("K8").Value = 1.033790650110609
$ws.Range("L8").Value = 1.033169602524252
$ws.Range("M8").Value = 1.040257290393607
$ws.Range("N8").Value = 1.014640665301424
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024932460479898
$ws.Range("D9").Value = 1.029681958638151
$ws.Range("E9").Value = 1.028668152687741
$ws.Range("F9").Value = 1.035048117974889
$ws.Range("I9").Value = 1.033675610330731
$ws.Range("J9").Value = 1.030801947246801
$ws.Range("K9").Value = 1.032865285369823
$ws.Range("L9").Value = 1.031854834870662
$ws.Range("M9").Value = 1.038213804044003
$ws.Range("N9").Value = 1.014289197700811
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023683282010027
$ws.Range("D10").Value = 1.028775586229752
$ws.Range("E10").Value = 1.027503448120229
$ws.Range("F10").Value = 1.033398427256921
$ws.Range("I10").Value = 1.03338849909931
$ws.Range("J10").Value = 1.030101306339908
$ws.Range("K10").Value = 1.032247196930741
$ws.Range("L10").Value = 1.030979681100072
$ws.Range("M10").Value = 1.036853346768265
$ws.Range("N10").Value = 1.014054545537573
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023143107527513
$ws.Range("D11").Value = 1.028383721799947
$ws.Range("E11").Value = 1.027000277782502
$ws.Range("F11").Value = 1.032685272322289
$ws.Range("I11").Value = 1.03326309183482
$ws.Range("J11").Value = 1.029797836577305
$ws.Range("K11").Value = 1.031979291728037
$ws.Range("L11").Value = 1.030601068622877
$ws.Range("M11").Value = 1.036264707779608
$ws.Range("N11").Value = 1.013952862918192
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022942573537348
$ws.Range("D12").Value = 1.028238257799269
$ws.Range("E12").Value = 1.026813553088556
$ws.Range("F12").Value = 1.032420551761554
$ws.Range("I12").Value = 1.033216347394399
$ws.Range("J12").Value = 1.029685102310476
$ws.Range("K12").Value = 1.031879740494633
$ws.Range("L12").Value = 1.030460486741249
$ws.Range("M12").Value = 1.03604612921009
$ws.Range("N12").Value = 1.013915082352928
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02298558371828
$ws.Range("D13").Value = 1.028269456136526
$ws.Range("E13").Value = 1.026853598214483
$ws.Range("F13").Value = 1.032477327200021
$ws.Range("I13").Value = 1.03322638158396
$ws.Range("J13").Value = 1.029709284722928
$ws.Range("K13").Value = 1.031901096333839
$ws.Range("L13").Value = 1.030490639673317
$ws.Range("M13").Value = 1.036093011957059
$ws.Range("N13").Value = 1.013923186908767
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023126529055142
$ws.Range("D14").Value = 1.028371695816187
$ws.Range("E14").Value = 1.02698483947755
$ws.Range("F14").Value = 1.032663386823786
$ws.Range("I14").Value = 1.033259231238326
$ws.Range("J14").Value = 1.029788518167925
$ws.Range("K14").Value = 1.031971063587772
$ws.Range("L14").Value = 1.030589447019301
$ws.Range("M14").Value = 1.036246638611454
$ws.Range("N14").Value = 1.013949740189226
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023213384863964
$ws.Range("D15").Value = 1.028434701283086
$ws.Range("E15").Value = 1.027065724815214
$ws.Range("F15").Value = 1.032778047769941
$ws.Range("I15").Value = 1.033279449466393
$ws.Range("J15").Value = 1.029837334939348
$ws.Range("K15").Value = 1.032014167545007
$ws.Range("L15").Value = 1.030650332370381
$ws.Range("M15").Value = 1.036341302121092
$ws.Range("N15").Value = 1.0139660990813
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02371914707698
$ws.Range("D16").Value = 1.028801605779038
$ws.Range("E16").Value = 1.027536866348834
$ws.Range("F16").Value = 1.033445781767021
$ws.Range("I16").Value = 1.033396799131767
$ws.Range("J16").Value = 1.030121444874581
$ws.Range("K16").Value = 1.032264971344635
$ws.Range("L16").Value = 1.031004815533029
$ws.Range("M16").Value = 1.036892422310734
$ws.Range("N16").Value = 1.014061292297318
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024036594314643
$ws.Range("D17").Value = 1.029031917281439
$ws.Range("E17").Value = 1.027832711282316
$ws.Range("F17").Value = 1.033864948046642
$ws.Range("I17").Value = 1.0334701190852
$ws.Range("J17").Value = 1.030299636833014
$ws.Range("K17").Value = 1.032422222728269
$ws.Range("L17").Value = 1.031227264057376
$ws.Range("M17").Value = 1.03723824598616
$ws.Range("N17").Value = 1.01412098424429
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02422182606303
$ws.Range("D18").Value = 1.029166311943822
$ws.Range("E18").Value = 1.028005383895875
$ws.Range("F18").Value = 1.034109553672503
$ws.Range("I18").Value = 1.033512780530591
$ws.Range("J18").Value = 1.03040356466049
$ws.Range("K18").Value = 1.032513918789534
$ws.Range("L18").Value = 1.031357046739499
$ws.Range("M18").Value = 1.037440002251967
$ws.Range("N18").Value = 1.014155794120797
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024284997155178
$ws.Range("D19").Value = 1.029212146814912
$ws.Range("E19").Value = 1.028064279615256
$ws.Range("F19").Value = 1.03419297701527
$ws.Range("I19").Value = 1.03352730917984
$ws.Range("J19").Value = 1.03043899988473
$ws.Range("K19").Value = 1.032545180332246
$ws.Range("L19").Value = 1.031401304704815
$ws.Range("M19").Value = 1.037508803220693
$ws.Range("N19").Value = 1.014167662117082
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024002527964142
$ws.Range("D20").Value = 1.029007201040406
$ws.Range("E20").Value = 1.027800958399197
$ws.Range("F20").Value = 1.033819963774723
$ws.Range("I20").Value = 1.033462263388221
$ws.Range("J20").Value = 1.030280519396467
$ws.Range("K20").Value = 1.032405353825619
$ws.Range("L20").Value = 1.031203394093473
$ws.Range("M20").Value = 1.037201137907478
$ws.Range("N20").Value = 1.014114580628485
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023085021094084
$ws.Range("D21").Value = 1.028341586228902
$ws.Range("E21").Value = 1.026946187347125
$ws.Range("F21").Value = 1.032608592010855
$ws.Range("I21").Value = 1.033249562317628
$ws.Range("J21").Value = 1.029765186218146
$ws.Range("K21").Value = 1.031950461051958
$ws.Range("L21").Value = 1.030560349282077
$ws.Range("M21").Value = 1.036201397502332
$ws.Range("N21").Value = 1.013941921213453
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022508790270477
$ws.Range("D22").Value = 1.027923620268236
$ws.Range("E22").Value = 1.02640977321608
$ws.Range("F22").Value = 1.031847978622917
$ws.Range("I22").Value = 1.033114888261232
$ws.Range("J22").Value = 1.029441106410095
$ws.Range("K22").Value = 1.031664224992191
$ws.Range("L22").Value = 1.030156340980658
$ws.Range("M22").Value = 1.035573215810761
$ws.Range("N22").Value = 1.013833299291644
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022814199766808
$ws.Range("D23").Value = 1.028145140914074
$ws.Range("E23").Value = 1.026694039768482
$ws.Range("F23").Value = 1.032251096810471
$ws.Range("I23").Value = 1.033186370472483
$ws.Range("J23").Value = 1.029612913464392
$ws.Range("K23").Value = 1.031815985293234
$ws.Range("L23").Value = 1.030370484590023
$ws.Range("M23").Value = 1.035906189103314
$ws.Range("N23").Value = 1.013890887766905
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024017920864999
$ws.Range("D24").Value = 1.029018369067263
$ws.Range("E24").Value = 1.027815305816707
$ws.Range("F24").Value = 1.033840289881677
$ws.Range("I24").Value = 1.033465813363556
$ws.Range("J24").Value = 1.030289157769618
$ws.Range("K24").Value = 1.032412976237111
$ws.Range("L24").Value = 1.031214179803295
$ws.Range("M24").Value = 1.037217905317601
$ws.Range("N24").Value = 1.014117474169822
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025417518495339
$ws.Range("D25").Value = 1.03003396552362
$ws.Range("E25").Value = 1.029120831238167
$ws.Range("F25").Value = 1.035688890714995
$ws.Range("I25").Value = 1.033785976363929
$ws.Range("J25").Value = 1.031073563227591
$ws.Range("K25").Value = 1.033104726761346
$ws.Range("L25").Value = 1.032194500267023
$ws.Range("M25").Value = 1.038741769224448
$ws.Range("N25").Value = 1.014380122279335
